$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record as row 3 (pushes old rows 3..28 down to 4..29) ---
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = 7
$ws.Range("B3").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C3").Value = 'Ñuble'
$ws.Range("D3").Value = 44817
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112037
$ws.Range("G3").Value = 'Cebollín'
$ws.Range("H3").Value = 'Sin especificar'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8500
$ws.Range("M3").Value = 8250
$ws.Range("N3").Value = '$/docena de atados'
$ws.Range("O3").Value = 'Provincia de Diguillín'
$ws.Range("P3").Value = 2750
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 'Hortaliza'

# --- Insert a second new record as row 22 (pushes the then-current rows 22..29 down to 23..30) ---
$ws.Rows("22:22").Insert()

$ws.Range("A22").Value = 7
$ws.Range("B22").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C22").Value = 'Ñuble'
$ws.Range("D22").Value = 44818
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112037
$ws.Range("G22").Value = 'Cebollín'
$ws.Range("H22").Value = 'Sin especificar'
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 8500
$ws.Range("M22").Value = 8250
$ws.Range("N22").Value = '$/docena de atados'
$ws.Range("O22").Value = 'Provincia de Diguillín'
$ws.Range("P22").Value = 2750
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 'Hortaliza'

Write-Output "Edit applied"
